$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 359.14285
$ws.Range("I33").Value = 219
$ws.Range("J33").Value = 546
$ws.Range("K33").Value = 219
$ws.Range("L33").Value = 546
$ws.Range("M33").Value = 10
$ws.Range("N33").Value = -1004
$ws.Range("H64").Value = 4160.8
$ws.Range("I64").Value = 4444.5
$ws.Range("J64").Value = 3971.6667
$ws.Range("K64").Value = 4444.5
$ws.Range("L64").Value = 3971.6667
$ws.Range("M64").Value = -4196.5
$ws.Range("N64").Value = -4467.6667
$ws.Range("H67").Value = 4160.8
$ws.Range("I67").Value = 4444.5
$ws.Range("J67").Value = 3971.6667
$ws.Range("K67").Value = 4444.5
$ws.Range("L67").Value = 3971.6667
$ws.Range("M67").Value = -3586.5
$ws.Range("N67").Value = -5687.6667
$ws.Range("H69").Value = 3388.25
$ws.Range("I69").Value = 2793
$ws.Range("J69").Value = 3586.6667
$ws.Range("K69").Value = 8379
$ws.Range("L69").Value = 10760.0001
$ws.Range("N69").Value = -12508.0001
$ws.Range("M69").Value = -7505
$ws.Range("H72").Value = 3388.25
$ws.Range("I72").Value = 2793
$ws.Range("J72").Value = 3586.6667
$ws.Range("K72").Value = 25137
$ws.Range("L72").Value = 32280.0003
$ws.Range("N72").Value = -41016.0003
$ws.Range("M72").Value = -20769
$ws.Range("H100").Value = 15625727
$ws.Range("I100").Value = 18519030
$ws.Range("J100").Value = 1884
$ws.Range("K100").Value = 18519030
$ws.Range("L100").Value = 1884
$ws.Range("M100").Value = -18518489
$ws.Range("N100").Value = -2966
$ws.Range("H103").Value = 1277
$ws.Range("I103").Value = 556.6667
$ws.Range("J103").Value = 1585.7142
$ws.Range("K103").Value = 1670.0001
$ws.Range("L103").Value = 4757.142599999999
$ws.Range("M103").Value = -1084.0001
$ws.Range("N103").Value = -5929.142599999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1140.3684
$ws.Range("I2").Value = 777.93335
$ws.Range("K2").Value = 777.93335
$ws.Range("M2").Value = -664.93335
$ws.Range("H32").Value = 4776.77
$ws.Range("I32").Value = 3805.0745
$ws.Range("K32").Value = 3805.0745
$ws.Range("M32").Value = -3518.0745
$ws.Range("H45").Value = 1252.375
$ws.Range("I45").Value = 1283.8
$ws.Range("J45").Value = 1200
$ws.Range("K45").Value = 1283.8
$ws.Range("L45").Value = 1200
$ws.Range("M45").Value = -906.8
$ws.Range("N45").Value = -1954
$ws.Range("H61").Value = 23810516
$ws.Range("I61").Value = 27027810
$ws.Range("J61").Value = 2535.6
$ws.Range("K61").Value = 27027810
$ws.Range("L61").Value = 2535.6
$ws.Range("M61").Value = -27027598
$ws.Range("N61").Value = -2959.6
$ws.Range("H63").Value = 18870544
$ws.Range("I63").Value = 2276.125
$ws.Range("J63").Value = 200005920
$ws.Range("K63").Value = 2276.125
$ws.Range("L63").Value = 200005920
$ws.Range("M63").Value = -1590.125
$ws.Range("N63").Value = -200007292
$ws.Range("H66").Value = 18870544
$ws.Range("I66").Value = 2276.125
$ws.Range("J66").Value = 200005920
$ws.Range("K66").Value = 11380.625
$ws.Range("L66").Value = 1000029600
$ws.Range("M66").Value = -7948.625
$ws.Range("N66").Value = -1000036464
$ws.Range("H74").Value = 2870.5715
$ws.Range("I74").Value = 1819.3334
$ws.Range("K74").Value = 1819.3334
$ws.Range("M74").Value = -945.3334
$ws.Range("H77").Value = 2870.5715
$ws.Range("I77").Value = 1819.3334
$ws.Range("K77").Value = 9096.666999999999
$ws.Range("M77").Value = -4728.666999999999
$ws.Range("H116").Value = 1140.3684
$ws.Range("I116").Value = 777.93335
$ws.Range("K116").Value = 777.93335
$ws.Range("M116").Value = 1516.06665
$ws.Range("H132").Value = 2243.5833
$ws.Range("I132").Value = 1565.907
$ws.Range("J132").Value = 3957.7058
$ws.Range("K132").Value = 4697.721
$ws.Range("L132").Value = 11873.1174
$ws.Range("M132").Value = -2167.721
$ws.Range("N132").Value = -16933.1174
$ws.Range("H136").Value = 23810516
$ws.Range("I136").Value = 27027810
$ws.Range("J136").Value = 2535.6
$ws.Range("K136").Value = 81083430
$ws.Range("L136").Value = 7606.799999999999
$ws.Range("M136").Value = -81080880
$ws.Range("N136").Value = -12706.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1140.3684
$ws.Range("I3").Value = 777.93335
$ws.Range("K3").Value = 777.93335
$ws.Range("M3").Value = -663.93335
$ws.Range("H105").Value = 47620264
$ws.Range("I105").Value = 52632690
$ws.Range("J105").Value = 2250
$ws.Range("K105").Value = 52632690
$ws.Range("L105").Value = 2250
$ws.Range("M105").Value = -52630943
$ws.Range("N105").Value = -5744
$ws.Range("H138").Value = 60874.285
$ws.Range("I138").Value = 90000
$ws.Range("J138").Value = 56020
$ws.Range("K138").Value = 90000
$ws.Range("L138").Value = 56020
$ws.Range("N138").Value = -66300
$ws.Range("M138").Value = -84860

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1179.9683
$ws.Range("I31").Value = 1039.2222
$ws.Range("K31").Value = 1039.2222
$ws.Range("M31").Value = -744.2221999999999
$ws.Range("H34").Value = 1179.9683
$ws.Range("I34").Value = 1039.2222
$ws.Range("K34").Value = 1039.2222
$ws.Range("M34").Value = -837.2221999999999
$ws.Range("H107").Value = 537.25
$ws.Range("I107").Value = 426.17648
$ws.Range("K107").Value = 426.17648
$ws.Range("M107").Value = 1493.82352

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 368.33334
$ws.Range("I46").Value = 105
$ws.Range("K46").Value = 315
$ws.Range("M46").Value = -224
$ws.Range("H98").Value = 290.25
$ws.Range("I98").Value = 164.75
$ws.Range("K98").Value = 494.25
$ws.Range("M98").Value = 1003.75
$ws.Range("H122").Value = 1605.3334
$ws.Range("I122").Value = 686.2222
$ws.Range("J122").Value = 2984
$ws.Range("K122").Value = 6175.999800000001
$ws.Range("L122").Value = 26856
$ws.Range("M122").Value = -3725.999800000001
$ws.Range("N122").Value = -31756

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 773.5
$ws.Range("I97").Value = 773.5
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 773.5
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -277.5
$ws.Range("N97").Value = $null
$ws.Range("H113").Value = 1926.6666
$ws.Range("I113").Value = 1980
$ws.Range("K113").Value = 1980
$ws.Range("M113").Value = 190
$ws.Range("H122").Value = 6000
$ws.Range("I122").Value = 6000
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 18000
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -15550
$ws.Range("N122").Value = $null
$ws.Range("H126").Value = 2134.9285
$ws.Range("I126").Value = 1642
$ws.Range("J126").Value = 2408.7778
$ws.Range("K126").Value = 4926
$ws.Range("L126").Value = 7226.3334
$ws.Range("M126").Value = -2456
$ws.Range("N126").Value = -12166.3334
$ws.Range("H132").Value = 2870.9062
$ws.Range("I132").Value = 2558.28
$ws.Range("K132").Value = 7674.84
$ws.Range("M132").Value = -5144.84

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 887.8889
$ws.Range("I22").Value = 831.8333
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 831.8333
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -536.8333
$ws.Range("N22").Value = -1590
$ws.Range("H27").Value = 887.8889
$ws.Range("I27").Value = 831.8333
$ws.Range("J27").Value = 1000
$ws.Range("K27").Value = 831.8333
$ws.Range("L27").Value = 1000
$ws.Range("M27").Value = -724.8333
$ws.Range("N27").Value = -1214
$ws.Range("H82").Value = 1978.1111
$ws.Range("I82").Value = 1907.0667
$ws.Range("K82").Value = 1907.0667
$ws.Range("M82").Value = -1546.0667
$ws.Range("H85").Value = 1978.1111
$ws.Range("I85").Value = 1907.0667
$ws.Range("K85").Value = 1907.0667
$ws.Range("M85").Value = -659.0667000000001
$ws.Range("H122").Value = 50001600
$ws.Range("I122").Value = 62501500
$ws.Range("K122").Value = 187504500
$ws.Range("M122").Value = -187502050

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 580.9655
$ws.Range("J107").Value = 737.82355
$ws.Range("L107").Value = 2213.47065
$ws.Range("N107").Value = -6053.470649999999
$ws.Range("H113").Value = 520.7059
$ws.Range("I113").Value = 351.76923
$ws.Range("J113").Value = 1055.30769
$ws.Range("K113").Value = 1055.30769
$ws.Range("M113").Value = 1114.69231
